# Economic Dashboard update - 2025-11-21
# Applies the weekly/daily data refresh described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: UI Initial Claims - new Lag3/Lag4 values populated ---
$ws.Range("T13").Value = 220000
$ws.Range("U13").Value = 232000

# --- Row 14: UI Continuing Claims - new Lag3/Lag4 values populated ---
$ws.Range("T14").Value = 1957000
$ws.Range("U14").Value = 1947000

# --- Row 29: 5yr, 5yr Forward (T5YIFR) - new date + rolling values ---
$ws.Range("N29").Value = 45981
$ws.Range("Q29").Value = 2.14
$ws.Range("R29").Value = 2.18
$ws.Range("S29").Value = 2.18
$ws.Range("T29").Value = 2.19

# --- Row 30: 10yr TIPS (T10YIE) - new date + rolling values ---
$ws.Range("N30").Value = 45981
$ws.Range("Q30").Value = 2.24
$ws.Range("R30").Value = 2.27
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = 2.28

# --- Row 39: Nominal Broad US Dollar Index - style-only change (remove highlight) ---
# Reuse an existing "no highlight" date-format cell so the workbook keeps the
# same shared style index instead of minting a new one.
$ws.Range("C39").Copy()
$ws.Range("N39").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N39").Value = 45975

# --- Row 42: Existing Home Sales - new period + rolling values ---
$ws.Range("C42").Value = 45931
$ws.Range("F42").Value = 4100000
$ws.Range("G42").Value = 4050000
$ws.Range("H42").Value = 4000000
$ws.Range("I42").Value = 4010000
$ws.Range("J42").Value = 3930000

# --- Row 43: Existing Home Sales Y/Y % Delta - new period + value ---
$ws.Range("C43").Value = 45931
$ws.Range("F43").Value = 0.0173697270471464

# --- Row 47: FFR (DFF) - new date ---
$ws.Range("N47").Value = 45980

# --- Row 48: 2y UST (DGS2) - new date + rolling values ---
$ws.Range("N48").Value = 45980
$ws.Range("R48").Value = 3.58
$ws.Range("S48").Value = 3.6
$ws.Range("U48").Value = $null

# --- Row 49: 5y UST (DGS5) - new date + rolling values ---
$ws.Range("N49").Value = 45980
$ws.Range("Q49").Value = 3.71
$ws.Range("R49").Value = 3.7
$ws.Range("S49").Value = 3.72
$ws.Range("U49").Value = $null

# --- Row 50: 10y UST (DGS10) - new date + rolling values ---
$ws.Range("N50").Value = 45980
$ws.Range("Q50").Value = 4.13
$ws.Range("R50").Value = 4.12
$ws.Range("S50").Value = 4.13
$ws.Range("U50").Value = $null

# --- Row 51: 30y Mortgage (MORTGAGE30US) - new date (+ highlight) + rolling values ---
# Reuse an existing "highlighted" date-format cell so the workbook keeps the
# same shared style index instead of minting a new one.
$ws.Range("N47").Copy()
$ws.Range("N51").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N51").Value = 45978
$ws.Range("Q51").Value = 6.26
$ws.Range("R51").Value = 6.24
$ws.Range("S51").Value = 6.22
$ws.Range("T51").Value = 6.17
$ws.Range("U51").Value = 6.19

# --- Row 52: BAA (DBAA) - new date + rolling values ---
$ws.Range("N52").Value = 45980
$ws.Range("Q52").Value = 5.92
$ws.Range("R52").Value = 5.91
$ws.Range("S52").Value = 5.9
$ws.Range("U52").Value = $null

Write-Output "Dashboard update applied"
